# docs/test-cases/autodiagnostico.xlsx
# Fill in the "Resultado Esperado" (H) and "Resultado Obtenido" (I) columns
# for the new "abrir modal" test cases (rows 7-10), plus the slightly
# different wording used on row 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("autodiagnostico")

$ws.Range("H7").Value = "El sistema debe permitir abrir el modal"
$ws.Range("I7").Value = "El modal se abre correctamente"

$ws.Range("H8").Value = "El sistema debe permitir abrir el modal"
$ws.Range("I8").Value = "El modal se abre correctamente"

$ws.Range("H9").Value = "El sistema debe permitir abrir el modal"
$ws.Range("I9").Value = "El modal se abre correctamente"

$ws.Range("H10").Value = "El sistema debe permitir abrir el modal"
$ws.Range("I10").Value = "El modal se abre correctamente"

$ws.Range("H11").Value = "El sistema debe permitir abrir el modal y mostrar la información"
$ws.Range("I11").Value = "El modal se abre y muestra la información correctamente"

# Reflect the saved selection state (the sheet was left scrolled/selected at E3)
$ws.Range("E3").Select() | Out-Null
